$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.791191816329956
$ws.Range("B1").Value = 1.957111239433289
$ws.Range("C1").Value = 2.137034177780151
$ws.Range("D1").Value = 2.130621910095215
$ws.Range("E1").Value = 2.983489274978638
